# Week 16 log + season sim from Week 17
# Applies updated cumulative season totals across YDS / OFF / DEF / ST / TURNS / PEN

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage samples to the running logs
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " -2 0 2 1 21 5 5 8 6 6 6 12 4 6 3 6 6 5 8 4 2 11 17 2"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 0 3 4 0 1 3 8 12 5 5 2 3 2 7 2 6 3 8 3 10 4 6 7 2 7 4 0 1 0 0 1 5 -2 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 18 9 5 -3 24 10 0 1 3 5 10 6 2 17 4 8 2 5 6 6 3 3 7 4 6 12"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 10 11 11 9 13 -1 5 18 18 5 7 3 14 8 50 3 8 9 4 6 9 16 22"

# ---------------------------------------------------------------------------
# OFF sheet: updated season cumulative offensive totals
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 172
$offWs.Range("E2").Value = 12
$offWs.Range("F2").Value = 53
$offWs.Range("G2").Value = 43
$offWs.Range("J2").Value = 21
$offWs.Range("N2").Value = 19
$offWs.Range("O2").Value = 17
$offWs.Range("P2").Value = 8

$offWs.Range("C3").Value = 176
$offWs.Range("E3").Value = 35
$offWs.Range("F3").Value = 114
$offWs.Range("G3").Value = 45
$offWs.Range("I3").Value = 65
$offWs.Range("J3").Value = 72
$offWs.Range("L3").Value = 273
$offWs.Range("M3").Value = 174
$offWs.Range("Q3").Value = 452

# ---------------------------------------------------------------------------
# DEF sheet: updated season cumulative defensive totals
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 203
$defWs.Range("D2").Value = 13
$defWs.Range("F2").Value = 63
$defWs.Range("G2").Value = 56
$defWs.Range("H2").Value = 5
$defWs.Range("I2").Value = 10
$defWs.Range("J2").Value = 30
$defWs.Range("N2").Value = 16
$defWs.Range("O2").Value = 19
$defWs.Range("P2").Value = 11

$defWs.Range("C3").Value = 154
$defWs.Range("D3").Value = 5
$defWs.Range("E3").Value = 31
$defWs.Range("F3").Value = 88
$defWs.Range("G3").Value = 34
$defWs.Range("H3").Value = 39
$defWs.Range("I3").Value = 53
$defWs.Range("J3").Value = 41
$defWs.Range("L3").Value = 245
$defWs.Range("M3").Value = 157
$defWs.Range("Q3").Value = 477

# ---------------------------------------------------------------------------
# ST sheet: updated season cumulative special-teams totals + per-game logs
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 76
$stWs.Range("D2").Value = 64
$stWs.Range("F2").Value = 242
$stWs.Range("G2").Value = 228
$stWs.Range("J2").Value = 106
$stWs.Range("K2").Value = 100

$stWs.Range("B3").Value = 39

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 60"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 20"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 19 29 13"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 58 63"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 18 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 20"

# ---------------------------------------------------------------------------
# TURNS sheet: updated season cumulative turnover totals
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B3").Value = 6
$turnsWs.Range("D3").Value = 8
$turnsWs.Range("E3").Value = 13

# ---------------------------------------------------------------------------
# PEN sheet: updated season cumulative penalty totals
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 19
$penWs.Range("D2").Value = 11
$penWs.Range("B3").Value = 14
